# ============================================================================
# Weekly CompStat data refresh (105th Precinct) — "New crime data collected"
#
# Updates the report header (volume/number + date range), then refreshes the
# Crime Complaints table (rows 15-21, 24-31) with the new week's figures.
# A handful of cells flip between "no data" placeholders (text "0" / "***.*")
# and real numbers (or vice versa) as categories go from empty to populated
# and back — those are handled with Range.Copy(Destination) from a stable
# same-styled reference cell so both the shared-string text AND the cell
# style/format travel together, instead of Value2 (which would silently
# coerce a numeric-looking string back into a number).
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header (rich-text cells): Volume 31 Number 16 -> 17, week dates ---
$ws.Range("A8").Value2 = "Volume 31   Number  17"
$ws.Range("C9").Value2 = "Report Covering the Week  4/22/2024  Through  4/28/2024"

# --- Crime Complaints table: plain numeric cell updates ---
$ws.Range("F15").Value2 = 2
$ws.Range("G15").Value2 = 1
$ws.Range("H15").Value2 = 100
$ws.Range("N15").Value2 = -35.714285714285
$ws.Range("D16").Value2 = 3
$ws.Range("E16").Value2 = -66.666666666666
$ws.Range("G16").Value2 = 14
$ws.Range("H16").Value2 = -57.142857142857
$ws.Range("I16").Value2 = 39
$ws.Range("J16").Value2 = 49
$ws.Range("K16").Value2 = -20.408163265306
$ws.Range("L16").Value2 = -26.415094339622
$ws.Range("M16").Value2 = -62.135922330097
$ws.Range("N16").Value2 = -87.888198757764
$ws.Range("C17").Value2 = 9
$ws.Range("D17").Value2 = 5
$ws.Range("E17").Value2 = 80
$ws.Range("F17").Value2 = 34
$ws.Range("G17").Value2 = 22
$ws.Range("H17").Value2 = 54.545454545454
$ws.Range("I17").Value2 = 144
$ws.Range("J17").Value2 = 127
$ws.Range("K17").Value2 = 13.385826771653
$ws.Range("L17").Value2 = 11.627906976744
$ws.Range("M17").Value2 = 51.578947368421
$ws.Range("N17").Value2 = 42.574257425742
$ws.Range("C18").Value2 = 4
$ws.Range("D18").Value2 = 11
$ws.Range("E18").Value2 = -63.636363636363
$ws.Range("F18").Value2 = 17
$ws.Range("G18").Value2 = 22
$ws.Range("H18").Value2 = -22.727272727272
$ws.Range("I18").Value2 = 59
$ws.Range("J18").Value2 = 75
$ws.Range("K18").Value2 = -21.333333333333
$ws.Range("L18").Value2 = 25.531914893617
$ws.Range("M18").Value2 = -48.695652173913
$ws.Range("N18").Value2 = -88.631984585741
$ws.Range("C19").Value2 = 11
$ws.Range("D19").Value2 = 8
$ws.Range("E19").Value2 = 37.5
$ws.Range("F19").Value2 = 44
$ws.Range("G19").Value2 = 44
$ws.Range("H19").Value2 = 0
$ws.Range("I19").Value2 = 208
$ws.Range("J19").Value2 = 188
$ws.Range("K19").Value2 = 10.63829787234
$ws.Range("L19").Value2 = 7.21649484536
$ws.Range("M19").Value2 = 37.74834437086
$ws.Range("N19").Value2 = 21.637426900584
$ws.Range("C20").Value2 = 7
$ws.Range("D20").Value2 = 4
$ws.Range("E20").Value2 = 75
$ws.Range("F20").Value2 = 36
$ws.Range("G20").Value2 = 33
$ws.Range("H20").Value2 = 9.090909090909
$ws.Range("I20").Value2 = 118
$ws.Range("J20").Value2 = 102
$ws.Range("K20").Value2 = 15.686274509803
$ws.Range("L20").Value2 = 49.367088607594
$ws.Range("M20").Value2 = 0
$ws.Range("N20").Value2 = -89.369369369369
$ws.Range("C21").Value2 = 32
$ws.Range("D21").Value2 = 31
$ws.Range("E21").Value2 = 3.225806451612
$ws.Range("F21").Value2 = 139
$ws.Range("G21").Value2 = 136
$ws.Range("H21").Value2 = 2.205882352941
$ws.Range("I21").Value2 = 579
$ws.Range("J21").Value2 = 547
$ws.Range("K21").Value2 = 5.850091407678
$ws.Range("L21").Value2 = 13.976377952755
$ws.Range("M21").Value2 = -2.525252525252
$ws.Range("N21").Value2 = -74.186357556843
$ws.Range("C24").Value2 = 20
$ws.Range("D24").Value2 = 20
$ws.Range("E24").Value2 = 0
$ws.Range("F24").Value2 = 88
$ws.Range("G24").Value2 = 81
$ws.Range("H24").Value2 = 8.641975308641
$ws.Range("I24").Value2 = 373
$ws.Range("J24").Value2 = 389
$ws.Range("K24").Value2 = -4.113110539845
$ws.Range("L24").Value2 = -17.66004415011
$ws.Range("M24").Value2 = 45.136186770428
$ws.Range("C25").Value2 = 3
$ws.Range("D25").Value2 = 5
$ws.Range("E25").Value2 = -40
$ws.Range("F25").Value2 = 16
$ws.Range("G25").Value2 = 16
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 72
$ws.Range("J25").Value2 = 66
$ws.Range("K25").Value2 = 9.090909090909
$ws.Range("L25").Value2 = -27.272727272727
$ws.Range("C26").Value2 = 14
$ws.Range("E26").Value2 = 7.692307692307
$ws.Range("F26").Value2 = 42
$ws.Range("G26").Value2 = 48
$ws.Range("H26").Value2 = -12.5
$ws.Range("I26").Value2 = 210
$ws.Range("J26").Value2 = 181
$ws.Range("K26").Value2 = 16.022099447513
$ws.Range("L26").Value2 = 6.598984771573
$ws.Range("M26").Value2 = 5
$ws.Range("F27").Value2 = 2
$ws.Range("G27").Value2 = 2
$ws.Range("H27").Value2 = 0
$ws.Range("L27").Value2 = 0
$ws.Range("C28").Value2 = 1
$ws.Range("E28").Value2 = 0
$ws.Range("F28").Value2 = 6
$ws.Range("G28").Value2 = 4
$ws.Range("H28").Value2 = 50
$ws.Range("I28").Value2 = 19
$ws.Range("J28").Value2 = 12
$ws.Range("K28").Value2 = 58.333333333333
$ws.Range("L28").Value2 = 46.153846153846
$ws.Range("G29").Value2 = 1
$ws.Range("H29").Value2 = 0
$ws.Range("G30").Value2 = 1
$ws.Range("H30").Value2 = 0

# --- C16: was empty ("0" placeholder, text), now has a real count (number) ---
# Set the value first, then pull the numeric-style format from a neighboring
# numeric cell in the same row so the style matches a real data cell (s=15)
# rather than leaving the old text style (s=14) behind.
$ws.Range("C16").Value2 = 1
$ws.Range("I16").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Rows 29-31: categories that went from populated back to empty this week.
# Copy-with-destination pulls both the shared "0" / "***.*" placeholder text
# AND the matching text-cell style (s=14) from stable reference cells
# (C22 = "0", E22 = "***.*") that are not touched elsewhere in this script.
$ws.Range("C22").Copy($ws.Range("D29"))
$ws.Range("E22").Copy($ws.Range("E29"))
$ws.Range("C22").Copy($ws.Range("D30"))
$ws.Range("E22").Copy($ws.Range("E30"))
$ws.Range("C22").Copy($ws.Range("G31"))
$ws.Range("E22").Copy($ws.Range("H31"))

# --- Column widths: best-fit narrowed for columns E (5) and H (8) now that
# the widest values in those columns are shorter than before. ---
$refWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $refWidth
$ws.Columns.Item(8).ColumnWidth = $refWidth
